$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear out the extra detail columns (D:O) on the remaining data row (row 2),
# keeping only No./Filename/Skills (columns A-C)
$ws.Range("D2:O2").Clear()

# Remove the other job rows (3-5), shifting everything up
$ws.Range("A3:A5").EntireRow.Delete()

# Turn the existing AutoFilter off, then reapply it over the now-smaller range
# (reapplying AutoFilter() on a range while AutoFilterMode is already on just
# toggles filtering off, so it must be disabled first)
$ws.AutoFilterMode = $false
$ws.Range("A1:O2").AutoFilter()

# Update the hidden _FilterDatabase defined name to match the new range
$nm = $wb.Names.Item("Sheet1!_FilterDatabase")
$nm.RefersTo = '=Sheet1!$A$1:$O$2'

# Reset the view/selection
$ws.Range("D8").Select()
